$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Actual" time column (E) for the team members rows with actual
# time spent on the third release.
$ws.Range("E29").Value = 14
$ws.Range("E30").Value = 11
$ws.Range("E31").Value = 10
$ws.Range("E32").Value = 11
$ws.Range("E33").Value = 9
$ws.Range("E34").Value = 11

$excel.CalculateFull()

# Restore the cursor/selection position as saved in the file.
$ws.Range("E37").Select()

$wb.Save()
